# Apply updated crypto price/volume data (GitHub Actions refresh)
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Cells.Item(2, 4).Value = "'70.722.54"
$ws.Cells.Item(2, 4).Style = "Normal"
$ws.Cells.Item(2, 5).Value = "'  -1.55%  "
$ws.Cells.Item(2, 5).Style = "Normal"

$ws.Cells.Item(3, 4).Value = "'2.540.26"
$ws.Cells.Item(3, 4).Style = "Normal"
$ws.Cells.Item(3, 5).Value = "'  -4.98%  "
$ws.Cells.Item(3, 5).Style = "Normal"

$ws.Cells.Item(4, 5).Value = "'  -0.03%  "
$ws.Cells.Item(4, 5).Style = "Normal"

$ws.Cells.Item(5, 4).Value = "'577.61"
$ws.Cells.Item(5, 4).Style = "Normal"
$ws.Cells.Item(5, 5).Value = "'  -3.42%  "
$ws.Cells.Item(5, 5).Style = "Normal"

$ws.Cells.Item(6, 4).Value = "'170.34"
$ws.Cells.Item(6, 4).Style = "Normal"
$ws.Cells.Item(6, 5).Value = "'  -2.71%  "
$ws.Cells.Item(6, 5).Style = "Normal"

$ws.Cells.Item(7, 5).Value = "'  +0.02%  "
$ws.Cells.Item(7, 5).Style = "Normal"

$ws.Cells.Item(8, 5).Value = "'  -2.69%  "
$ws.Cells.Item(8, 5).Style = "Normal"

$ws.Cells.Item(9, 5).Value = "'  -0.35%  "
$ws.Cells.Item(9, 5).Style = "Normal"

$ws.Cells.Item(10, 4).Value = "'2.538.03"
$ws.Cells.Item(10, 4).Style = "Normal"
$ws.Cells.Item(10, 5).Value = "'  -5.04%  "
$ws.Cells.Item(10, 5).Style = "Normal"

$ws.Cells.Item(11, 5).Value = "'  -0.30%  "
$ws.Cells.Item(11, 5).Style = "Normal"

$ws.Cells.Item(12, 5).Value = "'  -3.22%  "
$ws.Cells.Item(12, 5).Style = "Normal"

$ws.Cells.Item(13, 4).Value = "'4.81"
$ws.Cells.Item(13, 4).Style = "Normal"
$ws.Cells.Item(13, 5).Value = "'  -3.36%  "
$ws.Cells.Item(13, 5).Style = "Normal"

$ws.Cells.Item(14, 4).Value = "'3.005.65"
$ws.Cells.Item(14, 4).Style = "Normal"
$ws.Cells.Item(14, 5).Value = "'  -5.20%  "
$ws.Cells.Item(14, 5).Style = "Normal"

$ws.Cells.Item(15, 2).Value = "'ShibaInu"
$ws.Cells.Item(15, 2).Style = "Normal"
$ws.Cells.Item(15, 3).Value = "'https://coinranking.com/coin/xz24e0BjL+shibainu-shib"
$ws.Cells.Item(15, 3).Style = "Normal"
$ws.Cells.Item(15, 4).Value = "'0.0000181"
$ws.Cells.Item(15, 4).Style = "Normal"
$ws.Cells.Item(15, 5).Value = "'  -2.14%  "
$ws.Cells.Item(15, 5).Style = "Normal"

$ws.Cells.Item(16, 2).Value = "'WrappedBTC"
$ws.Cells.Item(16, 2).Style = "Normal"
$ws.Cells.Item(16, 3).Value = "'https://coinranking.com/coin/x4WXHge-vvFY+wrappedbtc-wbtc"
$ws.Cells.Item(16, 3).Style = "Normal"
$ws.Cells.Item(16, 4).Value = "'70.529.07"
$ws.Cells.Item(16, 4).Style = "Normal"
$ws.Cells.Item(16, 5).Value = "'  -1.71%  "
$ws.Cells.Item(16, 5).Style = "Normal"

$ws.Cells.Item(17, 4).Value = "'25.10"
$ws.Cells.Item(17, 4).Style = "Normal"
$ws.Cells.Item(17, 5).Value = "'  -4.16%  "
$ws.Cells.Item(17, 5).Style = "Normal"

$ws.Cells.Item(18, 4).Value = "'2.538.81"
$ws.Cells.Item(18, 4).Style = "Normal"
$ws.Cells.Item(18, 5).Value = "'  -4.61%  "
$ws.Cells.Item(18, 5).Style = "Normal"

$ws.Cells.Item(19, 4).Value = "'11.63"
$ws.Cells.Item(19, 4).Style = "Normal"
$ws.Cells.Item(19, 5).Value = "'  -4.80%  "
$ws.Cells.Item(19, 5).Style = "Normal"

$ws.Cells.Item(20, 4).Value = "'360.52"
$ws.Cells.Item(20, 4).Style = "Normal"
$ws.Cells.Item(20, 5).Value = "'  -2.78%  "
$ws.Cells.Item(20, 5).Style = "Normal"

$ws.Cells.Item(21, 4).Value = "'7.38"
$ws.Cells.Item(21, 4).Style = "Normal"
$ws.Cells.Item(21, 5).Value = "'  -10.64%  "
$ws.Cells.Item(21, 5).Style = "Normal"

$ws.Cells.Item(22, 5).Value = "'  -5.27%  "
$ws.Cells.Item(22, 5).Style = "Normal"

$ws.Cells.Item(23, 5).Value = "'  -1.33%  "
$ws.Cells.Item(23, 5).Style = "Normal"

$ws.Cells.Item(24, 5).Value = "'  -0.04%  "
$ws.Cells.Item(24, 5).Style = "Normal"

$ws.Cells.Item(25, 4).Value = "'69.59"
$ws.Cells.Item(25, 4).Style = "Normal"
$ws.Cells.Item(25, 5).Value = "'  -3.34%  "
$ws.Cells.Item(25, 5).Style = "Normal"

$ws.Cells.Item(26, 5).Value = "'  -5.64%  "
$ws.Cells.Item(26, 5).Style = "Normal"

$ws.Cells.Item(27, 5).Value = "'  -5.17%  "
$ws.Cells.Item(27, 5).Style = "Normal"

$ws.Cells.Item(28, 4).Value = "'2.669.57"
$ws.Cells.Item(28, 4).Style = "Normal"
$ws.Cells.Item(28, 5).Value = "'  -5.22%  "
$ws.Cells.Item(28, 5).Style = "Normal"

$ws.Cells.Item(29, 4).Value = "'0.999"
$ws.Cells.Item(29, 4).Style = "Normal"
$ws.Cells.Item(29, 5).Value = "'  -0.05%  "
$ws.Cells.Item(29, 5).Style = "Normal"

$ws.Cells.Item(30, 4).Value = "'0.0₃0919"
$ws.Cells.Item(30, 4).Style = "Normal"
$ws.Cells.Item(30, 5).Value = "'  -5.25%  "
$ws.Cells.Item(30, 5).Style = "Normal"

$ws.Cells.Item(31, 5).Value = "'  -2.12%  "
$ws.Cells.Item(31, 5).Style = "Normal"

$ws.Cells.Item(32, 4).Value = "'484.93"
$ws.Cells.Item(32, 4).Style = "Normal"
$ws.Cells.Item(32, 5).Value = "'  -3.38%  "
$ws.Cells.Item(32, 5).Style = "Normal"

$ws.Cells.Item(33, 5).Value = "'  -1.51%  "
$ws.Cells.Item(33, 5).Style = "Normal"

$ws.Cells.Item(34, 5).Value = "'  -3.32%  "
$ws.Cells.Item(34, 5).Style = "Normal"

$ws.Cells.Item(35, 5).Value = "'  +0.01%  "
$ws.Cells.Item(35, 5).Style = "Normal"

$ws.Cells.Item(36, 5).Value = "'  +5.87%  "
$ws.Cells.Item(36, 5).Style = "Normal"

$ws.Cells.Item(37, 4).Value = "'157.49"
$ws.Cells.Item(37, 4).Style = "Normal"
$ws.Cells.Item(37, 5).Value = "'  -2.95%  "
$ws.Cells.Item(37, 5).Style = "Normal"

$ws.Cells.Item(38, 4).Value = "'18.68"
$ws.Cells.Item(38, 4).Style = "Normal"
$ws.Cells.Item(38, 5).Value = "'  -4.43%  "
$ws.Cells.Item(38, 5).Style = "Normal"

$ws.Cells.Item(39, 4).Value = "'18.84"
$ws.Cells.Item(39, 4).Style = "Normal"
$ws.Cells.Item(39, 5).Value = "'  -1.23%  "
$ws.Cells.Item(39, 5).Style = "Normal"

$ws.Cells.Item(40, 4).Value = "'1.32"
$ws.Cells.Item(40, 4).Style = "Normal"
$ws.Cells.Item(40, 5).Value = "'  -4.55%  "
$ws.Cells.Item(40, 5).Style = "Normal"

$ws.Cells.Item(41, 5).Value = "'  +0.01%  "
$ws.Cells.Item(41, 5).Style = "Normal"

$ws.Cells.Item(42, 5).Value = "'  -5.61%  "
$ws.Cells.Item(42, 5).Style = "Normal"

$ws.Cells.Item(43, 4).Value = "'4.76"
$ws.Cells.Item(43, 4).Style = "Normal"
$ws.Cells.Item(43, 5).Value = "'  -4.81%  "
$ws.Cells.Item(43, 5).Style = "Normal"

$ws.Cells.Item(44, 2).Value = "'PolygonEcosystemToken"
$ws.Cells.Item(44, 2).Style = "Normal"
$ws.Cells.Item(44, 3).Value = "'https://coinranking.com/coin/iDZ0tG-wI+polygonecosystemtoken-pol"
$ws.Cells.Item(44, 3).Style = "Normal"
$ws.Cells.Item(44, 4).Value = "'0.319"
$ws.Cells.Item(44, 4).Style = "Normal"
$ws.Cells.Item(44, 5).Value = "'  -3.92%  "
$ws.Cells.Item(44, 5).Style = "Normal"

$ws.Cells.Item(45, 2).Value = "'dogwifhat"
$ws.Cells.Item(45, 2).Style = "Normal"
$ws.Cells.Item(45, 3).Value = "'https://coinranking.com/coin/sZUrmToWF+dogwifhat-wif"
$ws.Cells.Item(45, 3).Style = "Normal"
$ws.Cells.Item(45, 4).Value = "'2.44"
$ws.Cells.Item(45, 4).Style = "Normal"
$ws.Cells.Item(45, 5).Value = "'  -5.06%  "
$ws.Cells.Item(45, 5).Style = "Normal"

$ws.Cells.Item(46, 4).Value = "'38.40"
$ws.Cells.Item(46, 4).Style = "Normal"
$ws.Cells.Item(46, 5).Value = "'  -2.73%  "
$ws.Cells.Item(46, 5).Style = "Normal"

$ws.Cells.Item(47, 4).Value = "'144.70"
$ws.Cells.Item(47, 4).Style = "Normal"
$ws.Cells.Item(47, 5).Value = "'  -7.29%  "
$ws.Cells.Item(47, 5).Style = "Normal"

$ws.Cells.Item(48, 4).Value = "'3.54"
$ws.Cells.Item(48, 4).Style = "Normal"
$ws.Cells.Item(48, 5).Value = "'  -4.72%  "
$ws.Cells.Item(48, 5).Style = "Normal"

$ws.Cells.Item(49, 4).Value = "'0.528"
$ws.Cells.Item(49, 4).Style = "Normal"
$ws.Cells.Item(49, 5).Value = "'  -6.10%  "
$ws.Cells.Item(49, 5).Style = "Normal"

$ws.Cells.Item(50, 5).Value = "'  -6.49%  "
$ws.Cells.Item(50, 5).Style = "Normal"

$ws.Cells.Item(51, 4).Value = "'0.596"
$ws.Cells.Item(51, 4).Style = "Normal"
$ws.Cells.Item(51, 5).Value = "'  -1.61%  "
$ws.Cells.Item(51, 5).Style = "Normal"
